$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 42.75280866666666
$ws.Range("H2").Value = 128.258426
$ws.Range("I2").Value = 0.8529286054750734
$ws.Range("J2").Value = 0.8529286054750735
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.081040666666667
$ws.Range("N2").Value = 24.243122
$ws.Range("O2").Value = 0.4661250698616886
$ws.Range("P2").Value = 0.4661250698616886
$ws.Range("Q2").Value = 345.4871854495524
$ws.Range("R2").Value = 3109.384669045972
$ws.Range("S2").Value = 0.3975714058141012
$ws.Range("T2").Value = 0.3975714058141013

# Row 3
$ws.Range("G3").Value = 42.75280866666666
$ws.Range("H3").Value = 128.258426
$ws.Range("I3").Value = 0.8529286054750734
$ws.Range("J3").Value = 0.8529286054750735
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.000300666666668
$ws.Range("N3").Value = 21.000902
$ws.Range("O3").Value = 0.4037865631294714
$ws.Range("P3").Value = 0.4037865631294715
$ws.Range("Q3").Value = 299.2825150111391
$ws.Range("R3").Value = 2693.542635100252
$ws.Range("S3").Value = 0.3444011101995927
$ws.Range("T3").Value = 0.3444011101995929

# Row 4
$ws.Range("G4").Value = 42.75280866666666
$ws.Range("H4").Value = 128.258426
$ws.Range("I4").Value = 0.8529286054750734
$ws.Range("J4").Value = 0.8529286054750735
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.255294666666666
$ws.Range("N4").Value = 6.765884
$ws.Range("O4").Value = 0.1300883670088399
$ws.Range("P4").Value = 0.1300883670088399
$ws.Range("Q4").Value = 96.42018137095376
$ws.Range("R4").Value = 867.7816323385839
$ws.Range("S4").Value = 0.1109560894613794
$ws.Range("T4").Value = 0.1109560894613794

# Row 5
$ws.Range("G5").Value = 2.327094666666667
$ws.Range("H5").Value = 6.981284
$ws.Range("I5").Value = 0.04642608686423023
$ws.Range("J5").Value = 0.04642608686423023
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.081040666666667
$ws.Range("N5").Value = 24.243122
$ws.Range("O5").Value = 0.4661250698616886
$ws.Range("P5").Value = 0.4661250698616886
$ws.Range("Q5").Value = 18.80534663651645
$ws.Range("R5").Value = 169.248119728648
$ws.Range("S5").Value = 0.02164036298299414
$ws.Range("T5").Value = 0.02164036298299414

# Row 6
$ws.Range("G6").Value = 2.327094666666667
$ws.Range("H6").Value = 6.981284
$ws.Range("I6").Value = 0.04642608686423023
$ws.Range("J6").Value = 0.04642608686423023
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.000300666666668
$ws.Range("N6").Value = 21.000902
$ws.Range("O6").Value = 0.4037865631294714
$ws.Range("P6").Value = 0.4037865631294715
$ws.Range("Q6").Value = 16.29036234646312
$ws.Range("R6").Value = 146.613261118168
$ws.Range("S6").Value = 0.01874623005445782
$ws.Range("T6").Value = 0.01874623005445783

# Row 7
$ws.Range("G7").Value = 2.327094666666667
$ws.Range("H7").Value = 6.981284
$ws.Range("I7").Value = 0.04642608686423023
$ws.Range("J7").Value = 0.04642608686423023
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.255294666666666
$ws.Range("N7").Value = 6.765884
$ws.Range("O7").Value = 0.1300883670088399
$ws.Range("P7").Value = 0.1300883670088399
$ws.Range("Q7").Value = 5.248284190561778
$ws.Range("R7").Value = 47.234557715056
$ws.Range("S7").Value = 0.006039493826778264
$ws.Range("T7").Value = 0.006039493826778264

# Row 8
$ws.Range("G8").Value = 5.044817999999999
$ws.Range("H8").Value = 15.134454
$ws.Range("I8").Value = 0.1006453076606963
$ws.Range("J8").Value = 0.1006453076606963
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.081040666666667
$ws.Range("N8").Value = 24.243122
$ws.Range("O8").Value = 0.4661250698616886
$ws.Range("P8").Value = 0.4661250698616886
$ws.Range("Q8").Value = 40.76737941393199
$ws.Range("R8").Value = 366.9064147253879
$ws.Range("S8").Value = 0.0469133010645932
$ws.Range("T8").Value = 0.0469133010645932

# Row 9
$ws.Range("G9").Value = 5.044817999999999
$ws.Range("H9").Value = 15.134454
$ws.Range("I9").Value = 0.1006453076606963
$ws.Range("J9").Value = 0.1006453076606963
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.000300666666668
$ws.Range("N9").Value = 21.000902
$ws.Range("O9").Value = 0.4037865631294714
$ws.Range("P9").Value = 0.4037865631294715
$ws.Range("Q9").Value = 35.315242808612
$ws.Range("R9").Value = 317.837185277508
$ws.Range("S9").Value = 0.04063922287542081
$ws.Range("T9").Value = 0.04063922287542082

# Row 10
$ws.Range("G10").Value = 5.044817999999999
$ws.Range("H10").Value = 15.134454
$ws.Range("I10").Value = 0.1006453076606963
$ws.Range("J10").Value = 0.1006453076606963
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.255294666666666
$ws.Range("N10").Value = 6.765884
$ws.Range("O10").Value = 0.1300883670088399
$ws.Range("P10").Value = 0.1300883670088399
$ws.Range("Q10").Value = 11.377551129704
$ws.Range("R10").Value = 102.397960167336
$ws.Range("S10").Value = 0.01309278372068227
$ws.Range("T10").Value = 0.01309278372068227

